# Apply "estado" (column E) updates to Sheet1, flipping specific rows
# between "Trabaja" and "Nada" as described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows whose "estado" value becomes "Trabaja"
$trabajaRows = @(36,40,42,44,45,46,81,88,91,92,93,127,131,133,136,137,138,139,148,172,179,181,182,223,224,228,229,230,231,272,274,276,277,313,320,321,322,323,325,326,356,358)

# Rows whose "estado" value becomes "Nada"
$nadaRows = @(3,5,6,7,8,9,35,39,50,51,52,53,54,55,82,85,94,96,97,98,99,100,101,134,142,143,144,145,146,147,178,187,189,190,191,192,193,218,227,236,237,238,239,264,267,271,280,281,282,283,284,285,312,357,362,366,368)

foreach ($r in $trabajaRows) {
    $ws.Cells.Item($r, 5).Value = "Trabaja"
}

foreach ($r in $nadaRows) {
    $ws.Cells.Item($r, 5).Value = "Nada"
}
